$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value/type swaps (number <-> "NaN" text or corrected numbers) ---
$ws.Range("BV19").Value = "NaN"
$ws.Range("CM31").Value = "NaN"
$ws.Range("W36").Value = 1
$ws.Range("AY51").Value = "NaN"
$ws.Range("AY52").Value = 37
$ws.Range("AP87").Value = "NaN"
$ws.Range("AP88").Value = "NaN"
$ws.Range("L112").Value = 783
$ws.Range("L113").Value = 789
$ws.Range("J115").Value = 31
$ws.Range("J119").Value = 39
$ws.Range("J120").Value = 43
$ws.Range("J121").Value = 43
$ws.Range("J126").Value = 80
$ws.Range("J128").Value = 85
$ws.Range("W128").Value = "NaN"

# --- Minor numeric corrections in column J (rows 129-173, 176) ---
$ws.Range("J129").Value = 96
$ws.Range("J130").Value = 111
$ws.Range("J131").Value = 117
$ws.Range("J132").Value = 123
$ws.Range("J133").Value = 127
$ws.Range("J134").Value = 168
$ws.Range("J135").Value = 181
$ws.Range("J136").Value = 220
$ws.Range("J137").Value = 324
$ws.Range("J138").Value = 395
$ws.Range("J139").Value = 452
$ws.Range("J140").Value = 478
$ws.Range("J141").Value = 486
$ws.Range("J142").Value = 491
$ws.Range("J143").Value = 544
$ws.Range("J144").Value = 589
$ws.Range("J145").Value = 681
$ws.Range("J146").Value = 745
$ws.Range("J147").Value = 819
$ws.Range("J148").Value = 887
$ws.Range("J149").Value = 937
$ws.Range("J150").Value = 1119
$ws.Range("J151").Value = 1278
$ws.Range("J152").Value = 1517
$ws.Range("J153").Value = 1649
$ws.Range("J154").Value = 1841
$ws.Range("J155").Value = 1870
$ws.Range("J156").Value = 2045
$ws.Range("J157").Value = 2145
$ws.Range("J158").Value = 2497
$ws.Range("J159").Value = 2621
$ws.Range("J160").Value = 2796
$ws.Range("J161").Value = 2946
$ws.Range("J162").Value = 3070
$ws.Range("J163").Value = 3331
$ws.Range("J164").Value = 3514
$ws.Range("J165").Value = 3753
$ws.Range("J166").Value = 3925
$ws.Range("J167").Value = 4138
$ws.Range("J168").Value = 4347
$ws.Range("J169").Value = 4437
$ws.Range("J170").Value = 4612
$ws.Range("J171").Value = 4721
$ws.Range("J172").Value = 4841
$ws.Range("J173").Value = 5074
$ws.Range("J176").Value = 5347

# --- Add new row 178 with a full day of data ---
$ws.Range("A178").Value = 44072
$ws.Range("B178").Value = 599914
$ws.Range("C178").Value = 2696
$ws.Range("D178").Value = 78519
$ws.Range("E178").Value = 63872
$ws.Range("F178").Value = 207403
$ws.Range("G178").Value = 25012
$ws.Range("H178").Value = 3463
$ws.Range("I178").Value = 2804
$ws.Range("J178").Value = 5724
$ws.Range("K178").Value = 4983
$ws.Range("L178").Value = 10013
$ws.Range("M178").Value = 3704
$ws.Range("N178").Value = 19277
$ws.Range("O178").Value = 21972
$ws.Range("P178").Value = 4897
$ws.Range("Q178").Value = 4490
$ws.Range("R178").Value = 12109
$ws.Range("S178").Value = 7480
$ws.Range("T178").Value = 13814
$ws.Range("U178").Value = 11274
$ws.Range("V178").Value = 2779
$ws.Range("W178").Value = 1082
$ws.Range("X178").Value = 5758
$ws.Range("Y178").Value = 17261
$ws.Range("Z178").Value = 11506
$ws.Range("AA178").Value = 6730
$ws.Range("AB178").Value = 46086
$ws.Range("AC178").Value = 1008
$ws.Range("AD178").Value = 175
$ws.Range("AE178").Value = 260
$ws.Range("AF178").Value = 445
$ws.Range("AG178").Value = 93
$ws.Range("AH178").Value = 51
$ws.Range("AI178").Value = 241
$ws.Range("AJ178").Value = 1954
$ws.Range("AK178").Value = 2938
$ws.Range("AL178").Value = 35985
$ws.Range("AM178").Value = 6601
$ws.Range("AN178").Value = 2400
$ws.Range("AO178").Value = 36221
$ws.Range("AP178").Value = 907
$ws.Range("AQ178").Value = 20269
$ws.Range("AR178").Value = 1437
$ws.Range("AS178").Value = 7885
$ws.Range("AT178").Value = 1474
$ws.Range("AU178").Value = 1556
$ws.Range("AV178").Value = 4069
$ws.Range("AW178").Value = 1576
$ws.Range("AX178").Value = 932
$ws.Range("AY178").Value = 2464
$ws.Range("AZ178").Value = 2589
$ws.Range("BA178").Value = 45540
$ws.Range("BB178").Value = 11989
$ws.Range("BC178").Value = 2472
$ws.Range("BD178").Value = 7375
$ws.Range("BE178").Value = 3752
$ws.Range("BF178").Value = 278
$ws.Range("BG178").Value = 1393
$ws.Range("BH178").Value = 2580
$ws.Range("BI178").Value = 729
$ws.Range("BJ178").Value = 2002
$ws.Range("BK178").Value = 8179
$ws.Range("BL178").Value = 8250
$ws.Range("BM178").Value = 8190
$ws.Range("BN178").Value = 13749
$ws.Range("BO178").Value = 1866
$ws.Range("BP178").Value = 820
$ws.Range("BQ178").Value = 7054
$ws.Range("BR178").Value = 6247
$ws.Range("BS178").Value = 7382
$ws.Range("BT178").Value = 1472
$ws.Range("BU178").Value = 1471
$ws.Range("BV178").Value = 2792
$ws.Range("BW178").Value = 3081
$ws.Range("BX178").Value = 790
$ws.Range("BY178").Value = 4186
$ws.Range("BZ178").Value = 2396
$ws.Range("CA178").Value = 1277
$ws.Range("CB178").Value = 685
$ws.Range("CC178").Value = 1983
$ws.Range("CD178").Value = 1852
$ws.Range("CE178").Value = 1178
$ws.Range("CF178").Value = 917
$ws.Range("CG178").Value = 4708
$ws.Range("CH178").Value = 1329
$ws.Range("CI178").Value = 1138
$ws.Range("CJ178").Value = 1217
$ws.Range("CK178").Value = 1508
$ws.Range("CL178").Value = 1416
$ws.Range("CM178").Value = 1572
$ws.Range("CN178").Value = 1122
$ws.Range("CO178").Value = 1059
$ws.Range("CP178").Value = 1089
$ws.Range("CQ178").Value = 595
$ws.Range("CR178").Value = 2990
$ws.Range("CS178").Value = 958
$ws.Range("CT178").Value = 792
$ws.Range("CU178").Value = 731
$ws.Range("CV178").Value = 1267
$ws.Range("CW178").Value = 1137
$ws.Range("CX178").Value = 627
$ws.Range("CY178").Value = 720
$ws.Range("CZ178").Value = 884
$ws.Range("DA178").Value = 1157
$ws.Range("DB178").Value = 959
$ws.Range("DC178").Value = 1097
$ws.Range("DD178").Value = 869
$ws.Range("DE178").Value = 312
$ws.Range("DF178").Value = 331
$ws.Range("DG178").Value = 664
$ws.Range("DH178").Value = 578
$ws.Range("DI178").Value = 399
$ws.Range("DJ178").Value = 530
$ws.Range("DK178").Value = 324
$ws.Range("DL178").Value = 586
$ws.Range("DM178").Value = 699
$ws.Range("DN178").Value = 508
$ws.Range("DO178").Value = 475
$ws.Range("DP178").Value = 358
$ws.Range("DQ178").Value = 512
$ws.Range("DR178").Value = 115903
$ws.Range("DS178").Value = 253259
$ws.Range("DT178").Value = 9917
$ws.Range("DU178").Value = 109292
$ws.Range("DV178").Value = 68587
$ws.Range("DW178").Value = 28584
$ws.Range("DX178").Value = 8596

# --- Update active selection to match the saved view ---
$ws.Range("DR177").Select()

